$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits at the end
#    of the first paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append a brand-new paragraph after the last paragraph in the body,
#    containing the new scaffolding/BootStrap instructions, and re-create
#    the "_GoBack" bookmark around its very end (mirroring what Word does
#    after the most recent edit location).
$endRange = $d.Range($d.Content.End, $d.Content.End)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:firstLine="720"/></w:pPr>' +
    '<w:r><w:t>Wri</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">te an error validation message in the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>message.properties</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> file. </w:t></w:r>' +
    '<w:r><w:t>Add some sample data to your BootStrap.groovy?</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$endRange.InsertXML($newParagraphXml) | Out-Null
